$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 331, shifting existing rows 331:435 down to 332:436.
$ws.Range("A331").EntireRow.Insert()

# Populate the newly inserted row 331 with the new data record.
$ws.Range("A331").Value = 3
$ws.Range("B331").Value = "Femacal de La Calera"
$ws.Range("C331").Value = "Coquimbo"
$ws.Range("D331").Value = 44876
$ws.Range("E331").Value = 5
$ws.Range("F331").Value = 100114013
$ws.Range("G331").Value = "Zanahoria"
$ws.Range("H331").Value = "Sin especificar"
$ws.Range("I331").Value = "Primera"
$ws.Range("J331").Value = 230
$ws.Range("K331").Value = 12000
$ws.Range("L331").Value = 12500
$ws.Range("M331").Value = 12261
$ws.Range("N331").Value = "$/saco 20 kilos"
$ws.Range("O331").Value = "Provincia de Quillota"
$ws.Range("P331").Value = 613
$ws.Range("Q331").Value = 20
$ws.Range("R331").Value = "Hortaliza"
